# Updated symbol list on Mon Dec 12 20:45:31 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking quotes as TEXT (the sheet's
# original XML uses inlineStr for every cell). Assigning a plain numeric
# string via Range.Value lets the engine auto-coerce it to a real number,
# which would change the cell's stored type. To keep these as text - exactly
# like the source data - force text entry by setting NumberFormat to "@"
# before the write, then clear the (now unneeded) explicit format back off
# the cell so no stray style index is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Column D ("Price") updates
Set-TextValue "D2"  "274.45"
Set-TextValue "D4"  "6.262"
Set-TextValue "D5"  "0.06186"
Set-TextValue "D6"  "3.573"
Set-TextValue "D7"  "1.525"
Set-TextValue "D8"  "6.534"
Set-TextValue "D9"  "0.8225"
Set-TextValue "D10" "0.1649"
Set-TextValue "D11" "0.08295"
Set-TextValue "D12" "0.03471"
Set-TextValue "D15" "3.767"
Set-TextValue "D16" "0.001628"
Set-TextValue "D17" "0.04680"
Set-TextValue "D18" "0.006261"
Set-TextValue "D19" "0.006133"
Set-TextValue "D23" "2.310"
Set-TextValue "D24" "0.01388"
Set-TextValue "D28" "0.0002737"
Set-TextValue "D40" "0.04725"
Set-TextValue "D41" "0.005300"
Set-TextValue "D42" "0.007045"
Set-TextValue "D43" "0.1103"
Set-TextValue "D44" "0.01147"
Set-TextValue "D45" "0.00006058"
Set-TextValue "D47" "0.7230"

# Column E ("Volume(1h)") label swaps - plain text, no numeric coercion risk
$ws.Range("E41").Value = "40CEJICEJI"
$ws.Range("E42").Value = "41KickTokenKICKBestin24h"

Write-Host "Applied cryptos price/volume refresh."
